# Updated local file list spreadsheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "AllFiles" value that used to live in B2 (this also drops the
# now-unused "AllFiles" entry from the shared strings table on save).
$ws.Range("B2").ClearContents()

# Restore the last-used selection recorded in the sheet view.
$ws.Range("D13").Select() | Out-Null
